# New crime data collected - weekly update for 121st Precinct CompStat report
# Moves the reporting week forward one week (2/10-2/16/2025 -> 2/17-2/23/2025,
# volume "Number" 7 -> 8) and refreshes the crime-complaint statistics table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: bump the bulletin "Number" and the reporting week's two dates.
# These are rich-text cells; Characters() edits only the targeted substring.
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 1).Text = "8"
$ws.Range("C9").Characters(27, 9).Text = "2/17/2025"
$ws.Range("C9").Characters(47, 9).Text = "2/23/2025"

# ---------------------------------------------------------------------------
# Helper donor cells used to carry over the correct number format when a
# cell's underlying type flips between the "blank" text placeholder
# (style of C14, values "0" / "***.*") and an actual numeric entry
# (style of F15 for whole numbers, L15 for one-decimal percentages).
# PasteSpecial (formats only) reassigns the cell style without touching
# any other cell's value.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

function Set-AsBlankText {
    param($cellRef, $text)
    $ws.Range("C14").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($cellRef).Value = "'" + $text
}

function Set-AsWholeNumber {
    param($cellRef, $value)
    $ws.Range("F15").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($cellRef).Value = $value
}

function Set-AsPercentNumber {
    param($cellRef, $value)
    $ws.Range("L15").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($cellRef).Value = $value
}

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-AsWholeNumber "C15" 1
$ws.Range("I15").Value = 4
$ws.Range("L15").Value = 100

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 26.666666666666
$ws.Range("I17").Value = 29
$ws.Range("J17").Value = 24
$ws.Range("K17").Value = 20.833333333333
$ws.Range("L17").Value = -14.705882352941

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 1
Set-AsWholeNumber "D18" 3
Set-AsPercentNumber "E18" -66.666666666666
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = 63.636363636363
$ws.Range("L18").Value = 200

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -43.333333333333
$ws.Range("I19").Value = 41
$ws.Range("J19").Value = 57
$ws.Range("K19").Value = -28.070175438596
$ws.Range("L19").Value = 0

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-AsBlankText "G20" "0"
Set-AsBlankText "H20" "***.*"
$ws.Range("I20").Value = 18
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 50

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = 6.896551724137
$ws.Range("I21").Value = 118
$ws.Range("J21").Value = 105
$ws.Range("K21").Value = 12.380952380952
$ws.Range("L21").Value = 10.280373831775

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
Set-AsBlankText "C23" "0"

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -3.448275862068
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 125
$ws.Range("H24").Value = -3.2
$ws.Range("I24").Value = 234
$ws.Range("J24").Value = 239
$ws.Range("K24").Value = -2.092050209205
$ws.Range("L24").Value = 27.173913043478

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -19.047619047619
$ws.Range("F25").Value = 79
$ws.Range("H25").Value = -14.130434782608
$ws.Range("I25").Value = 159
$ws.Range("J25").Value = 158
$ws.Range("K25").Value = 0.632911392405
$ws.Range("L25").Value = 34.745762711864

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 5
$ws.Range("E26").Value = -61.538461538461
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = -40.476190476190
$ws.Range("I26").Value = 55
$ws.Range("J26").Value = 79
$ws.Range("K26").Value = -30.379746835443
$ws.Range("L26").Value = -25.675675675675

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
Set-AsWholeNumber "C27" 1
Set-AsWholeNumber "D27" 1
Set-AsPercentNumber "E27" 0
Set-AsWholeNumber "G27" 1
Set-AsPercentNumber "H27" 100
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = -20

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-AsBlankText "C28" "0"
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 10
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = 66.666666666666
$ws.Range("L28").Value = 25

# ---------------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------------
Set-AsBlankText "D31" "0"
Set-AsBlankText "E31" "***.*"

$excel.CutCopyMode = 0
